$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 1246
$ws.Range("L3").Value = 1256
$ws.Range("I4").Value = 1828
$ws.Range("J4").Value = 1857
$ws.Range("K4").Value = 1745
$ws.Range("L4").Value = 347
$ws.Range("L5").Value = 83
$ws.Range("L6").Value = 1253
$ws.Range("I7").Value = 26293
$ws.Range("J7").Value = 29329
$ws.Range("K7").Value = 27538
$ws.Range("L7").Value = 4185

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L4").Value = 14
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 129
$ws.Range("L8").Value = 257
$ws.Range("L10").Value = 26
$ws.Range("L15").Value = 30
$ws.Range("L19").Value = 124
$ws.Range("I20").Value = 641
$ws.Range("K20").Value = 672
$ws.Range("L20").Value = 113
$ws.Range("L23").Value = 41
$ws.Range("L27").Value = 47
$ws.Range("L29").Value = 201
$ws.Range("L31").Value = 44
$ws.Range("L33").Value = 185
$ws.Range("L37").Value = 143
$ws.Range("L42").Value = 134
$ws.Range("L46").Value = 10
$ws.Range("L51").Value = 55
$ws.Range("L52").Value = 83
$ws.Range("L54").Value = 92
$ws.Range("L60").Value = 24
$ws.Range("I63").Value = 249
$ws.Range("J63").Value = 206
$ws.Range("L63").Value = 15
$ws.Range("L64").Value = 29
$ws.Range("L65").Value = 84
$ws.Range("L67").Value = 154
$ws.Range("L72").Value = 18
$ws.Range("L73").Value = 34
$ws.Range("L79").Value = 113
$ws.Range("L85").Value = 213
$ws.Range("L86").Value = 32
$ws.Range("L88").Value = 62
$ws.Range("L89").Value = 54
$ws.Range("L91").Value = 55
$ws.Range("L92").Value = 10
$ws.Range("L95").Value = 65
$ws.Range("L98").Value = 31
$ws.Range("I101").Value = 26293
$ws.Range("J101").Value = 29329
$ws.Range("K101").Value = 27538
$ws.Range("L101").Value = 4185

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 33
$ws.Range("L3").Value = 46
$ws.Range("L7").Value = 129

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 57
$ws.Range("L3").Value = 91
$ws.Range("L4").Value = 18
$ws.Range("L7").Value = 213

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L3").Value = 22
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 83

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 74
$ws.Range("L3").Value = 80
$ws.Range("L7").Value = 257

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 39
$ws.Range("L3").Value = 68
$ws.Range("L7").Value = 185

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L3").Value = 22
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 44
$ws.Range("L6").Value = 44
$ws.Range("L7").Value = 143

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 84

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L3").Value = 46
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 154

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L6").Value = 46
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 69
$ws.Range("L3").Value = 68
$ws.Range("L5").Value = 3
$ws.Range("L6").Value = 53
$ws.Range("L7").Value = 201

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 37
$ws.Range("L7").Value = 124

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L3").Value = 32
$ws.Range("L6").Value = 55
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 10

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L2").Value = 24
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L4").Value = 5
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L3").Value = 31
$ws.Range("I4").Value = 43
$ws.Range("K4").Value = 32
$ws.Range("L6").Value = 36
$ws.Range("I7").Value = 641
$ws.Range("K7").Value = 672
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("L2").Value = 13
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 10

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L2").Value = 12
$ws.Range("L3").Value = 21
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 47

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L2").Value = 5
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L3").Value = 18
$ws.Range("L4").Value = 3
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("L2").Value = 9
$ws.Range("L7").Value = 18

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 14
